$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.085.13"
$ws.Range("E2").Value = "  +0.24%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.631.01"
$ws.Range("E3").Value = "  -0.78%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.12%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.42"
$ws.Range("E5").Value = "  -0.49%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.55%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -1.52%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0621"
$ws.Range("E9").Value = "  -2.78%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.65"
$ws.Range("E10").Value = "  -4.83%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.05%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.857.62"
$ws.Range("E12").Value = "  -0.79%  "

# Row 13 - now Polkadot (was WrappedEther)
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.19"
$ws.Range("E13").Value = "  -1.56%  "

# Row 14 - now WrappedEther (was Polkadot)
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.618.37"
$ws.Range("E14").Value = "  -1.00%  "

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.529"
$ws.Range("E15").Value = "  -2.72%  "

# Row 16 - WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.068.68"
$ws.Range("E16").Value = "  +0.02%  "

# Row 17 - ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0745"
$ws.Range("E17").Value = "  -2.11%  "

# Row 18 - Litecoin
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.67"
$ws.Range("E18").Value = "  -2.77%  "

# Row 19 - Dai
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.00"
$ws.Range("E19").Value = "  -0.02%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.56"
$ws.Range("E20").Value = "  -0.51%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.27"
$ws.Range("E21").Value = "  -2.10%  "

# Row 22 - Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.58"
$ws.Range("E22").Value = "  -3.29%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  -2.03%  "

# Row 24 - Stellar
$ws.Range("E24").Value = "  +1.67%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.33"
$ws.Range("E25").Value = "  +0.22%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.06%  "

# Row 27 - Toncoin
$ws.Range("E27").Value = "  -3.71%  "

# Row 28 - Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.75"
$ws.Range("E28").Value = "  -2.05%  "

# Row 29 - EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.33"
$ws.Range("E29").Value = "  -1.13%  "

# Row 30 - PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("E30").Value = "  -0.75%  "

# Row 31 - Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0484"
$ws.Range("E31").Value = "  -2.35%  "

# Row 32 - Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.14"
$ws.Range("E32").Value = "  -3.87%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -4.55%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -2.40%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -1.74%  "

# Row 36 - Maker
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.131.15"
$ws.Range("E36").Value = "  +0.08%  "

# Row 37 - ARBITRUM
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.855"
$ws.Range("E37").Value = "  -5.49%  "

# Row 38 - MXToken
$ws.Range("E38").Value = "  -1.15%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -3.33%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -2.03%  "

# Row 41 - Quant
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.43"
$ws.Range("E41").Value = "  -0.45%  "

# Row 42 - RocketPoolETH
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.767.04"
$ws.Range("E42").Value = "  -0.83%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  -4.62%  "

# Row 44 - FraxShare
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.14"
$ws.Range("E44").Value = "  -5.61%  "

# Row 45 - BabyDogeCoin
$ws.Range("E45").Value = "  -1.78%  "

# Row 46 - Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "54.65"
$ws.Range("E46").Value = "  -3.24%  "

# Row 47 - Cronos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0524"
$ws.Range("E47").Value = "  +0.52%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  -0.13%  "

# Row 49 - Mantle
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.414"
$ws.Range("E49").Value = "  -0.10%  "

# Row 50 - EnergySwap
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.52"
$ws.Range("E50").Value = "  -3.24%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  +0.18%  "
